$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 110.42857
$ws.Range("I11").Value = 110.42857
$ws.Range("K11").Value = 110.42857
$ws.Range("M11").Value = 29.57143000000001
$ws.Range("H33").Value = 146.55556
$ws.Range("J33").Value = 233.33333
$ws.Range("L33").Value = 233.33333
$ws.Range("N33").Value = -691.3333299999999
$ws.Range("H38").Value = 2373.2
$ws.Range("I38").Value = 413.6
$ws.Range("J38").Value = 4332.8
$ws.Range("K38").Value = 1240.8
$ws.Range("L38").Value = 12998.4
$ws.Range("M38").Value = -868.8000000000002
$ws.Range("N38").Value = -13742.4
$ws.Range("H39").Value = 472.1111
$ws.Range("I39").Value = 270.2
$ws.Range("K39").Value = 810.5999999999999
$ws.Range("M39").Value = -514.5999999999999
$ws.Range("H43").Value = 3250
$ws.Range("I43").Value = 3250
$ws.Range("K43").Value = 3250
$ws.Range("M43").Value = -3181
$ws.Range("H92").Value = 226.625
$ws.Range("I92").Value = 189.57143
$ws.Range("K92").Value = 189.57143
$ws.Range("M92").Value = 1058.42857
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H137").Value = 1934.15
$ws.Range("I137").Value = 1077.9
$ws.Range("J137").Value = 2790.4
$ws.Range("K137").Value = 3233.7
$ws.Range("L137").Value = 8371.200000000001
$ws.Range("M137").Value = -683.7000000000003
$ws.Range("N137").Value = -13471.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4002534.8
$ws.Range("I32").Value = 974.9048
$ws.Range("K32").Value = 974.9048
$ws.Range("M32").Value = -687.9048
$ws.Range("H45").Value = 2334.4
$ws.Range("I45").Value = 1491.6666
$ws.Range("J45").Value = 3598.5
$ws.Range("K45").Value = 1491.6666
$ws.Range("L45").Value = 3598.5
$ws.Range("M45").Value = -1114.6666
$ws.Range("N45").Value = -4352.5
$ws.Range("H61").Value = 1171.35
$ws.Range("I61").Value = 1171.35
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1171.35
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -959.3499999999999
$ws.Range("N61").ClearContents()
$ws.Range("H110").Value = 2918
$ws.Range("I110").Value = 2054.8333
$ws.Range("J110").Value = 5507.5
$ws.Range("K110").Value = 2054.8333
$ws.Range("L110").Value = 5507.5
$ws.Range("M110").Value = -9.833299999999781
$ws.Range("N110").Value = -9597.5
$ws.Range("H122").Value = 2175.8572
$ws.Range("I122").Value = 1746.5454
$ws.Range("K122").Value = 5239.6362
$ws.Range("M122").Value = -2789.6362
$ws.Range("H132").Value = 1776.7858
$ws.Range("I132").Value = 1792.9487
$ws.Range("K132").Value = 5378.8461
$ws.Range("M132").Value = -2848.8461
$ws.Range("H136").Value = 1171.35
$ws.Range("I136").Value = 1171.35
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3514.05
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -964.0499999999997
$ws.Range("N136").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 431.33334
$ws.Range("I22").Value = 416.8
$ws.Range("K22").Value = 416.8
$ws.Range("M22").Value = -243.8
$ws.Range("H97").Value = 4875.8335
$ws.Range("I97").Value = 4451
$ws.Range("J97").Value = 7000
$ws.Range("K97").Value = 4451
$ws.Range("L97").Value = 7000
$ws.Range("M97").Value = -3460
$ws.Range("N97").Value = -8982
$ws.Range("H99").Value = 2610
$ws.Range("I99").Value = 2610
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2610
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1112
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 1923.5714
$ws.Range("I105").Value = 1900
$ws.Range("K105").Value = 1900
$ws.Range("M105").Value = -153
$ws.Range("H107").Value = 4824.25
$ws.Range("I107").Value = 798.6667
$ws.Range("K107").Value = 798.6667
$ws.Range("M107").Value = 1121.3333
$ws.Range("H134").Value = 1344.2307
$ws.Range("I134").Value = 1170.04
$ws.Range("K134").Value = 3510.12
$ws.Range("M134").Value = -975.1199999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4880.2
$ws.Range("J31").Value = 6100
$ws.Range("L31").Value = 6100
$ws.Range("N31").Value = -6690
$ws.Range("H34").Value = 4880.2
$ws.Range("J34").Value = 6100
$ws.Range("L34").Value = 6100
$ws.Range("N34").Value = -6504
$ws.Range("H122").Value = 1090.8572
$ws.Range("I122").Value = 1022.6667
$ws.Range("K122").Value = 3068.0001
$ws.Range("M122").Value = -618.0001000000002
$ws.Range("H134").Value = 1470.4117
$ws.Range("I134").Value = 1531.4286
$ws.Range("K134").Value = 4594.2858
$ws.Range("M134").Value = -2059.2858

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 284127
$ws.Range("I4").Value = 284127
$ws.Range("K4").Value = 852381
$ws.Range("M4").Value = -852269
$ws.Range("H11").Value = 39138.06
$ws.Range("I11").Value = 94323.08
$ws.Range("J11").Value = 4975.905
$ws.Range("K11").Value = 282969.24
$ws.Range("L11").Value = 14927.715
$ws.Range("M11").Value = -282829.24
$ws.Range("N11").Value = -15207.715
$ws.Range("H33").Value = 62.909092
$ws.Range("I33").Value = 56
$ws.Range("J33").Value = 75
$ws.Range("K33").Value = 336
$ws.Range("L33").Value = 450
$ws.Range("M33").Value = -53
$ws.Range("N33").Value = -1016
$ws.Range("H103").Value = 2206.2727
$ws.Range("J103").Value = 2206.2727
$ws.Range("L103").Value = 6618.8181
$ws.Range("N103").Value = -8376.8181
$ws.Range("H117").Value = 89
$ws.Range("I117").Value = 60
$ws.Range("J117").Value = 118
$ws.Range("K117").Value = 180
$ws.Range("L117").Value = 354
$ws.Range("M117").Value = 3262
$ws.Range("N117").Value = -7238
$ws.Range("H122").Value = 728.44446
$ws.Range("I122").Value = 654
$ws.Range("K122").Value = 5886
$ws.Range("M122").Value = -3436
$ws.Range("H124").Value = 1400
$ws.Range("I124").Value = 1400
$ws.Range("K124").Value = 4200
$ws.Range("M124").Value = 710
$ws.Range("H140").Value = 2099.5454
$ws.Range("I140").Value = 1784.5
$ws.Range("K140").Value = 5353.5
$ws.Range("M140").Value = -173.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 13238714
$ws.Range("I11").Value = 12500812
$ws.Range("J11").Value = 15600000
$ws.Range("K11").Value = 12500812
$ws.Range("L11").Value = 15600000
$ws.Range("M11").Value = -12500673
$ws.Range("N11").Value = -15600278
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H107").Value = 189.88889
$ws.Range("I107").Value = 169.83333
$ws.Range("J107").Value = 230
$ws.Range("K107").Value = 169.83333
$ws.Range("L107").Value = 230
$ws.Range("M107").Value = 1750.16667
$ws.Range("N107").Value = -4070
$ws.Range("H122").Value = 1847.4706
$ws.Range("I122").Value = 1814.8572
$ws.Range("K122").Value = 5444.571599999999
$ws.Range("M122").Value = -2994.571599999999
$ws.Range("H132").Value = 36434.4
$ws.Range("I132").Value = 47499.316
$ws.Range("K132").Value = 142497.948
$ws.Range("M132").Value = -139967.948

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2708.7727
$ws.Range("I46").Value = 1625.8572
$ws.Range("J46").Value = 3214.1333
$ws.Range("K46").Value = 1625.8572
$ws.Range("L46").Value = 3214.1333
$ws.Range("M46").Value = -1437.8572
$ws.Range("N46").Value = -3590.1333
$ws.Range("H100").Value = 7416.6665
$ws.Range("I100").Value = 4750
$ws.Range("K100").Value = 4750
$ws.Range("M100").Value = -4209
$ws.Range("H122").Value = 3860
$ws.Range("I122").Value = 3899.5
$ws.Range("K122").Value = 11698.5
$ws.Range("M122").Value = -9248.5
$ws.Range("H132").Value = 8765.556
$ws.Range("I132").Value = 6269.2856
$ws.Range("K132").Value = 18807.8568
$ws.Range("M132").Value = -16277.8568
$ws.Range("H140").Value = 101856
$ws.Range("J140").Value = 69141.336
$ws.Range("L140").Value = 69141.336
$ws.Range("N140").Value = -79501.336

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5168.625
$ws.Range("I122").Value = 2499.6667
$ws.Range("K122").Value = 7499.000100000001
$ws.Range("M122").Value = -5049.000100000001
$ws.Range("H132").Value = 1571.25
$ws.Range("I132").Value = 1596
$ws.Range("K132").Value = 4788
$ws.Range("M132").Value = -2258

